$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9 (for "fr-observation-multiresistant-microorganism-document"),
# which pushes all subsequent rows down by one.
$ws.Range("A9:K9").Insert()

# Copy style from the row above into the freshly inserted row, row 9 (only used columns).
$ws.Range("A8:K8").Copy()
$ws.Range("A9:K9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A9").Value = "fr-observation-multiresistant-microorganism-document"
$ws.Range("B9").Value = "Observation - FR Observation Multiresistant Microorganisms Identification Document"
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = "TerminologieCISIS - Terminologie des concepts non trouvés dans les autres terminologies#MED-144"
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = "dateTimeĵ, Periodĵ, Timingĵ, instantĵ"
$ws.Range("H9").Value = "stringĵ"
$ws.Range("I9").Value = "optional"
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = ""

# Row 11 ("fr-observation-pregnancy-history-document"): Time Types value changes.
$ws.Range("G11").Value = "Periodĵ"
